# Update automàtic: dades i banners [2026-02-06 04:49]
# Applies the diff between the previous meteocat extraction snapshot and the
# new one: updated DATA_EXTRACCIO timestamps and several measurement columns
# (HUMITAT_MITJANA_DIA, PRESSIO_ATMOSFERICA, RADIACIO_GLOBAL, RATXA_VENT_MAX,
# TEMPERATURA_MINIMA_DIA, TEMPERATURA_MITJANA_DIA) for rows 2-36.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Direct text/date/measurement updates (safe from numeric auto-conversion) ---
$ws.Range("E2").Value = '2026-02-06 04:47:57'
$ws.Range("N2").Value = '-2.6 °C 4:19 TU'
$ws.Range("O2").Value = '-1.5 °C'
$ws.Range("E3").Value = '2026-02-06 04:48:00'
$ws.Range("E4").Value = '2026-02-06 04:48:02'
$ws.Range("J4").Value = '992.8 hPa'
$ws.Range("N4").Value = '10.3 °C 4:22 TU'
$ws.Range("O4").Value = '12.9 °C'
$ws.Range("E5").Value = '2026-02-06 04:48:05'
$ws.Range("J5").Value = '993.4 hPa'
$ws.Range("N5").Value = '6.3 °C 4:28 TU'
$ws.Range("O5").Value = '8.3 °C'
$ws.Range("E6").Value = '2026-02-06 04:48:07'
$ws.Range("J6").Value = '994.4 hPa'
$ws.Range("N6").Value = '13.8 °C 4:07 TU'
$ws.Range("E7").Value = '2026-02-06 04:48:10'
$ws.Range("J7").Value = '994.2 hPa'
$ws.Range("L7").Value = '41.8 km/h - 251º 4:06 TU'
$ws.Range("N7").Value = '9.5 °C 4:29 TU'
$ws.Range("E8").Value = '2026-02-06 04:48:13'
$ws.Range("N8").Value = '4.5 °C 4:07 TU'
$ws.Range("O8").Value = '6.1 °C'
$ws.Range("E9").Value = '2026-02-06 04:48:15'
$ws.Range("N9").Value = '1.2 °C 4:22 TU'
$ws.Range("O9").Value = '2.4 °C'
$ws.Range("E10").Value = '2026-02-06 04:48:18'
$ws.Range("N10").Value = '3.6 °C 4:03 TU'
$ws.Range("O10").Value = '5.2 °C'
$ws.Range("E11").Value = '2026-02-06 04:48:21'
$ws.Range("J11").Value = '995.1 hPa'
$ws.Range("N11").Value = '4.0 °C 4:29 TU'
$ws.Range("O11").Value = '4.7 °C'
$ws.Range("E12").Value = '2026-02-06 04:48:23'
$ws.Range("N12").Value = '10.0 °C 4:29 TU'
$ws.Range("O12").Value = '12.6 °C'
$ws.Range("E13").Value = '2026-02-06 04:48:25'
$ws.Range("N13").Value = '4.8 °C 4:22 TU'
$ws.Range("O13").Value = '6.7 °C'
$ws.Range("E14").Value = '2026-02-06 04:48:28'
$ws.Range("E15").Value = '2026-02-06 04:48:31'
$ws.Range("J15").Value = '993.4 hPa'
$ws.Range("N15").Value = '3.7 °C 4:18 TU'
$ws.Range("O15").Value = '7.3 °C'
$ws.Range("E16").Value = '2026-02-06 04:48:33'
$ws.Range("N16").Value = '3.4 °C 4:01 TU'
$ws.Range("O16").Value = '4.2 °C'
$ws.Range("E17").Value = '2026-02-06 04:48:36'
$ws.Range("J17").Value = '996.7 hPa'
$ws.Range("N17").Value = '1.8 °C 4:21 TU'
$ws.Range("O17").Value = '3.2 °C'
$ws.Range("E18").Value = '2026-02-06 04:48:39'
$ws.Range("N18").Value = '-5.3 °C 4:29 TU'
$ws.Range("O18").Value = '-4.8 °C'
$ws.Range("E19").Value = '2026-02-06 04:48:41'
$ws.Range("J19").Value = '997.0 hPa'
$ws.Range("E20").Value = '2026-02-06 04:48:44'
$ws.Range("O20").Value = '-2.0 °C'
$ws.Range("E21").Value = '2026-02-06 04:48:46'
$ws.Range("J21").Value = '994.3 hPa'
$ws.Range("O21").Value = '5.1 °C'
$ws.Range("E22").Value = '2026-02-06 04:48:49'
$ws.Range("N22").Value = '4.6 °C 4:29 TU'
$ws.Range("O22").Value = '8.5 °C'
$ws.Range("E23").Value = '2026-02-06 04:48:52'
$ws.Range("J23").Value = '993.6 hPa'
$ws.Range("E24").Value = '2026-02-06 04:48:54'
$ws.Range("J24").Value = '992.4 hPa'
$ws.Range("E25").Value = '2026-02-06 04:48:57'
$ws.Range("J25").Value = '995.7 hPa'
$ws.Range("N25").Value = '1.1 °C 4:29 TU'
$ws.Range("O25").Value = '2.2 °C'
$ws.Range("E26").Value = '2026-02-06 04:48:59'
$ws.Range("N26").Value = '-1.1 °C 4:00 TU'
$ws.Range("O26").Value = '-0.4 °C'
$ws.Range("E27").Value = '2026-02-06 04:49:02'
$ws.Range("J27").Value = '993.2 hPa'
$ws.Range("N27").Value = '5.9 °C 4:28 TU'
$ws.Range("O27").Value = '7.6 °C'
$ws.Range("E28").Value = '2026-02-06 04:49:05'
$ws.Range("J28").Value = '996.2 hPa'
$ws.Range("O28").Value = '3.1 °C'
$ws.Range("E29").Value = '2026-02-06 04:49:08'
$ws.Range("N29").Value = '7.1 °C 4:17 TU'
$ws.Range("O29").Value = '12.0 °C'
$ws.Range("E30").Value = '2026-02-06 04:49:10'
$ws.Range("E31").Value = '2026-02-06 04:49:13'
$ws.Range("J31").Value = '996.6 hPa'
$ws.Range("N31").Value = '4.1 °C 4:29 TU'
$ws.Range("O31").Value = '5.1 °C'
$ws.Range("E32").Value = '2026-02-06 04:49:16'
$ws.Range("J32").Value = '994.8 hPa'
$ws.Range("K32").Value = '-0.1 MJ/m2'
$ws.Range("N32").Value = '12.8 °C 4:29 TU'
$ws.Range("O32").Value = '15.2 °C'
$ws.Range("E33").Value = '2026-02-06 04:49:18'
$ws.Range("N33").Value = '5.2 °C 4:26 TU'
$ws.Range("O33").Value = '6.9 °C'
$ws.Range("E34").Value = '2026-02-06 04:49:21'
$ws.Range("N34").Value = '4.5 °C 4:29 TU'
$ws.Range("O34").Value = '8.4 °C'
$ws.Range("E35").Value = '2026-02-06 04:49:23'
$ws.Range("N35").Value = '-3.4 °C 4:01 TU'
$ws.Range("E36").Value = '2026-02-06 04:49:26'
$ws.Range("J36").Value = '996.3 hPa'
$ws.Range("K36").Value = '-0.1 MJ/m2'
$ws.Range("N36").Value = '9.4 °C 4:22 TU'
$ws.Range("O36").Value = '11.8 °C'

# --- Percentage-looking values: must be forced as literal text, not converted to numbers ---
$scratch = $ws.Range("ZZ500")
$scratch.NumberFormat = "@"
$scratch.Value = '75%'
$scratch.Copy()
$ws.Range("H3").PasteSpecial(-4163)
$scratch.Value = '58%'
$scratch.Copy()
$ws.Range("H4").PasteSpecial(-4163)
$scratch.Value = '75%'
$scratch.Copy()
$ws.Range("H5").PasteSpecial(-4163)
$scratch.Value = '92%'
$scratch.Copy()
$ws.Range("H8").PasteSpecial(-4163)
$scratch.Value = '83%'
$scratch.Copy()
$ws.Range("H15").PasteSpecial(-4163)
$scratch.Value = '98%'
$scratch.Copy()
$ws.Range("H19").PasteSpecial(-4163)
$scratch.Value = '71%'
$scratch.Copy()
$ws.Range("H20").PasteSpecial(-4163)
$scratch.Value = '85%'
$scratch.Copy()
$ws.Range("H21").PasteSpecial(-4163)
$scratch.Value = '81%'
$scratch.Copy()
$ws.Range("H22").PasteSpecial(-4163)
$scratch.Value = '93%'
$scratch.Copy()
$ws.Range("H23").PasteSpecial(-4163)
$scratch.Value = '82%'
$scratch.Copy()
$ws.Range("H26").PasteSpecial(-4163)
$scratch.Value = '97%'
$scratch.Copy()
$ws.Range("H27").PasteSpecial(-4163)
$scratch.Value = '61%'
$scratch.Copy()
$ws.Range("H29").PasteSpecial(-4163)
$scratch.Value = '76%'
$scratch.Copy()
$ws.Range("H30").PasteSpecial(-4163)
$scratch.Value = '76%'
$scratch.Copy()
$ws.Range("H34").PasteSpecial(-4163)
$scratch.Value = '95%'
$scratch.Copy()
$ws.Range("H35").PasteSpecial(-4163)
$scratch.Value = '66%'
$scratch.Copy()
$ws.Range("H36").PasteSpecial(-4163)
$scratch.Clear()
$excel.CutCopyMode = $false

